# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.127881588408715, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 1, 1.094976487407548)
    3 = @(1.459612070389937, 1.667794583268128, 3.900430680208489, 8.660232485948974, 1, 15.68806981981553)
    4 = @(0.003994804209775715, 3099.503889238888, 3.900430680208489, 645.3272768299601, 1, 3748.735591553267)
    5 = @(3.230985683306322, 10.29869402782916, 3.900430680208489, 8.660232485948974, 1, 26.09034287729295)
    6 = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 0, 12.59312877619104)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("F$row").Value = $vals[4]
    $ws.Range("G$row").Value = $vals[5]
}
